# TestData.xlsx update: refresh SignIn / CreateAccount sample rows,
# drop the duplicate second data row on each sheet, and re-point the
# mailto hyperlinks at the new e-mail address.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "SignIn"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SignIn")

# Drop existing hyperlinks (they'll be rebuilt below) and remove the
# redundant third row (row 2 is the only data row that survives).
$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(3).Delete()

# Refresh the row-2 values (email changes, the rest stays the same).
$ws1.Range("A2").Value = "Y"
$ws1.Range("B2").Value = "SignIn"
$ws1.Range("C2").Value = "testjaga277171@gmail.com"
$ws1.Range("D2").Value = "jaga@12345"

# Re-create the hyperlinks and make sure the cells keep the Hyperlink look.
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:testjaga277171@gmail.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:jaga@12345") | Out-Null
$ws1.Range("C2").Style = "Hyperlink"
$ws1.Range("D2").Style = "Hyperlink"

# ---------------------------------------------------------------
# Sheet "CreateAccount"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CreateAccount")

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Range("A2").Value = "Y"
$ws2.Range("B2").Value = "CreateAccount"
$ws2.Range("C2").Value = "Jaga"
$ws2.Range("D2").Value = "Waran"
$ws2.Range("E2").Value = "08/30/1991"
$ws2.Range("F2").Value = "testjaga277171@gmail.com"
$ws2.Range("G2").Value = "IT Company"
$ws2.Range("H2").Value = "Sholinganallur"
# I2 (zip code, "600119") is unchanged - leave its value/format untouched.
$ws2.Range("J2").Value = "Chennai"
$ws2.Range("K2").Value = "TamilNadu"
$ws2.Range("L2").Value = "India"

# M2 (phone number) takes on the same "quote-prefixed" text look as I2
# (zip code) -- copy its format over before writing the new value.
$ws2.Range("I2").Copy() | Out-Null
$ws2.Range("M2").PasteSpecial(-4122) | Out-Null
$ws2.Range("M2").Value = "'9876543212"

$ws2.Range("N2").Value = "jaga@12345"
$ws2.Range("O2").Value = "jaga@12345"

$ws2.Hyperlinks.Add($ws2.Range("F2"), "mailto:testjaga277171@gmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("N2"), "mailto:jaga@12345") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("O2"), "mailto:jaga@12345") | Out-Null
$ws2.Range("F2").Style = "Hyperlink"
$ws2.Range("N2").Style = "Hyperlink"
$ws2.Range("O2").Style = "Hyperlink"

# ---------------------------------------------------------------
# Selection bookkeeping (CreateAccount first so SignIn ends up the
# active/front tab, matching the saved view state).
# ---------------------------------------------------------------
$ws2.Range("F7").Select() | Out-Null
$ws1.Range("D18").Select() | Out-Null
